$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.084.66'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = '''1.929.77'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.24%  '
$ws.Range("D4").Value = '''1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''326.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.36%  '
$ws.Range("D6").Value = '''1.004'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("D7").Value = '''0.4609'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.99%  '
$ws.Range("D8").Value = '''0.3832'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.04%  '
$ws.Range("D9").Value = '''0.07755'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("D10").Value = '''0.9814'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.32%  '
$ws.Range("D11").Value = '''22.56'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.91%  '
$ws.Range("D12").Value = '''1.951.02'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.76%  '
$ws.Range("D13").Value = '''6.983'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.42%  '
$ws.Range("D14").Value = '''5.697'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").Value = '''0.07056'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").Value = '''1.006'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").Value = '''84.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.10%  '
$ws.Range("D18").Value = '''0.000009562'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.89%  '
$ws.Range("D19").Value = '''16.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.70%  '
$ws.Range("D20").Value = '''1.003'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("D21").Value = '''29.097.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = '''5.348'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.26%  '
$ws.Range("E23").Value = '  +1.27%  '
$ws.Range("D24").Value = '''2.079'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.87%  '
$ws.Range("D25").Value = '''157.66'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.61%  '
$ws.Range("D26").Value = '''19.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.10%  '
$ws.Range("D27").Value = '''5.684'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("D28").Value = '''118.08'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.19%  '
$ws.Range("D29").Value = '''1.853'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.18%  '
$ws.Range("D30").Value = '''0.09344'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.22%  '
$ws.Range("D31").Value = '''0.8662'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.32%  '
$ws.Range("D32").Value = '''5.129'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.51%  '
$ws.Range("D33").Value = '''1.249'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("D34").Value = '''3.023'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.97%  '
$ws.Range("D35").Value = '''0.05712'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.72%  '
$ws.Range("E36").Value = '  +1.13%  '
$ws.Range("D37").Value = '''1.004'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("D38").Value = '''0.02048'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.06%  '
$ws.Range("D39").Value = '''3.072'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +14.23%  '
$ws.Range("D40").Value = '''7.550'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.69%  '
$ws.Range("D41").Value = '''0.5520'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.04%  '
$ws.Range("D42").Value = '''0.1755'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '''9.382'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.38%  '
$ws.Range("B44").Value = 'PEPE'
$ws.Range("C44").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D44").Value = '''0.000002879'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.49%  '
$ws.Range("D45").Value = '''2.220'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.83%  '
$ws.Range("D46").Value = '''0.5204'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.98%  '
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("D48").Value = '''0.06926'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.39%  '
$ws.Range("D49").Value = '''1.780'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.82%  '
$ws.Range("D50").Value = '''110.46'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.15%  '
$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").Value = '''1.003'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.13%  '
